$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Review audit CV Maju Tecnology" row (row 5); rows below shift up.
$ws.Rows.Item(5).Delete()

# Update the week number from 38 to 41 for all data rows (now rows 2-9).
$ws.Range("B2:B9").Value = 41

# Shorten the CS Tuparev follow-up task text (now on row 8 after the deletion).
$ws.Range("C8").Value = "Follow up kasus CS Tuparev & CS Patrol"

# Update the selected cell to match the saved view state.
$ws.Range("F16").Select()
